# "support blank row and colum start"
# Sheet3's data block is shifted so it no longer starts at column A / row 1:
# the old A1:E12 block of "same" values becomes a B2:F12 block (column A and
# row 1 are cleared), a new column F is filled in alongside B:E, and two
# more rows (13 and 14) of data are appended below, introducing three new
# shared strings (NewF14 / NewG14 / NewG13). Sheet3 also becomes the active
# (selected) sheet, with the cursor left on G14; Sheet1 loses its previous
# tabSelected flag as a result.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet3")

# Capture the existing shared "same" text before we clear column A.
$sameValue = $ws.Range("B2").Value()

# Blank out the old leading column (A) and leading row (1) so the used
# range starts at B2 instead of A1.
$ws.Range("A1:A12").ClearContents()
$ws.Range("A1:E1").ClearContents()

# Extend the existing B:E block rightwards into column F for rows 2-12.
$ws.Range("F2:F12").Value = $sameValue

# Row 13 repeats the same value across B:F, plus a brand new G13 string.
$ws.Range("B13:F13").Value = $sameValue
$ws.Range("F14").Value = "NewF14"
$ws.Range("G14").Value = "NewG14"
$ws.Range("G13").Value = "NewG13"

# Sheet3 becomes the active sheet/tab, with G14 as the selected cell.
$ws.Activate()
$ws.Range("G14").Select()
